$d = $word.ActiveDocument

# Update the date header
$d.Content.Find.Execute("2025-05-29 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-30 Friday", 2)

# Update the table of division problems/answers
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "96÷3=32, 0"
$t.Cell(1, 2).Range.Text = "23÷8=2, 7"
$t.Cell(1, 3).Range.Text = "74÷5=14, 4"
$t.Cell(1, 4).Range.Text = "58÷5=11, 3"
$t.Cell(1, 5).Range.Text = "95÷5=19, 0"

$t.Cell(5, 1).Range.Text = "13÷9=1, 4"
$t.Cell(5, 2).Range.Text = "94÷8=11, 6"
$t.Cell(5, 3).Range.Text = "70÷8=8, 6"
$t.Cell(5, 4).Range.Text = "37÷7=5, 2"
$t.Cell(5, 5).Range.Text = "95÷3=31, 2"

$t.Cell(9, 1).Range.Text = "43÷9=4, 7"
$t.Cell(9, 2).Range.Text = "24÷9=2, 6"
$t.Cell(9, 3).Range.Text = "19÷2=9, 1"
$t.Cell(9, 4).Range.Text = "95÷2=47, 1"
$t.Cell(9, 5).Range.Text = "58÷8=7, 2"

$t.Cell(13, 1).Range.Text = "85÷5=17, 0"
$t.Cell(13, 2).Range.Text = "87÷8=10, 7"
$t.Cell(13, 3).Range.Text = "17÷8=2, 1"
$t.Cell(13, 4).Range.Text = "79÷2=39, 1"
$t.Cell(13, 5).Range.Text = "84÷7=12, 0"

$t.Cell(17, 1).Range.Text = "43÷9=4, 7"
$t.Cell(17, 2).Range.Text = "54÷6=9, 0"
$t.Cell(17, 3).Range.Text = "58÷3=19, 1"
$t.Cell(17, 4).Range.Text = "25÷3=8, 1"
$t.Cell(17, 5).Range.Text = "29÷9=3, 2"
